$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 447, shifting existing rows 447:530 down to 450:533.
# xlShiftDown = -4121
$ws.Range("A447:R449").Insert(-4121)

# Fill in the new rows (447:449) with the new weekly price report for
# "Cuatro cascos verde" dated 2021-11-04 (serial 44504).

# Row 447: Cuatro cascos verde - Primera
$ws.Range("A447").Value = 2
$ws.Range("B447").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C447").Value = "Coquimbo"
$ws.Range("D447").Value = 44504
$ws.Range("E447").Value = 4
$ws.Range("F447").Value = 100112002
$ws.Range("G447").Value = "Pimiento"
$ws.Range("H447").Value = "Cuatro cascos verde"
$ws.Range("I447").Value = "Primera"
$ws.Range("J447").Value = 1000
$ws.Range("K447").Value = 23000
$ws.Range("L447").Value = 24000
$ws.Range("M447").Value = 23500
$ws.Range("N447").Value = "`$/caja 18 kilos"
$ws.Range("O447").Value = "Provincia de Limarí"
$ws.Range("P447").Value = 1306
$ws.Range("Q447").Value = 18
$ws.Range("R447").Value = "Hortaliza"

# Row 448: Cuatro cascos verde - Segunda
$ws.Range("A448").Value = 2
$ws.Range("B448").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C448").Value = "Coquimbo"
$ws.Range("D448").Value = 44504
$ws.Range("E448").Value = 4
$ws.Range("F448").Value = 100112002
$ws.Range("G448").Value = "Pimiento"
$ws.Range("H448").Value = "Cuatro cascos verde"
$ws.Range("I448").Value = "Segunda"
$ws.Range("J448").Value = 700
$ws.Range("K448").Value = 18000
$ws.Range("L448").Value = 19000
$ws.Range("M448").Value = 18500
$ws.Range("N448").Value = "`$/caja 18 kilos"
$ws.Range("O448").Value = "Provincia de Limarí"
$ws.Range("P448").Value = 1028
$ws.Range("Q448").Value = 18
$ws.Range("R448").Value = "Hortaliza"

# Row 449: Cuatro cascos verde - Tercera
$ws.Range("A449").Value = 2
$ws.Range("B449").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C449").Value = "Coquimbo"
$ws.Range("D449").Value = 44504
$ws.Range("E449").Value = 4
$ws.Range("F449").Value = 100112002
$ws.Range("G449").Value = "Pimiento"
$ws.Range("H449").Value = "Cuatro cascos verde"
$ws.Range("I449").Value = "Tercera"
$ws.Range("J449").Value = 500
$ws.Range("K449").Value = 13000
$ws.Range("L449").Value = 14000
$ws.Range("M449").Value = 13500
$ws.Range("N449").Value = "`$/caja 18 kilos"
$ws.Range("O449").Value = "Provincia de Limarí"
$ws.Range("P449").Value = 750
$ws.Range("Q449").Value = 18
$ws.Range("R449").Value = "Hortaliza"
